# Apply updated values to Sheet1 for rows 2-25, columns B,C,D,F,G,I,L,M,O
$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets.Item("Sheet1")

$data = @(
    @(1.030328928315214, 0.3333604406678603, 0.02747446327848735, 0.6667722330233588, 0.002420739447516543, 0.6297596269737085, 0.2776079781274561, 0.2395579635942653, 2.269282291836191),
    @(0.9208930879644868, 0.3155287654023766, 0.02556991596924973, 0.6669166347670767, 0.002423500990466341, 0.6401161304543166, 0.2766578962475705, 0.2230850627558638, 2.284275099860949),
    @(0.8536182700579502, 0.3045284538312956, 0.02439392387783101, 0.6675541946790347, 0.002425287591582055, 0.6469492478080099, 0.276233648720023, 0.21301787775176, 2.295328758516376),
    @(0.8261846477008703, 0.3000330892038505, 0.02391307122421438, 0.667951826492569, 0.002426038597582614, 0.64985286684324, 0.2761008018259474, 0.2089275411436802, 2.300297154256697),
    @(0.8216282439325369, 0.2992858830085652, 0.02383312870833265, 0.6680261710192639, 0.002426164689731807, 0.6503421970440577, 0.2760811617137051, 0.2082490832722712, 2.30115015116246),
    @(0.8532483636001871, 0.3044678785504118, 0.02438744547468019, 0.6675589995205442, 0.002425297626939031, 0.6469879251522705, 0.2762316949519317, 0.2129626646319167, 2.295393886574701),
    @(0.9926132097600089, 0.3272229753136742, 0.02681916180105759, 0.6667079748748392, 0.002421672781090517, 0.6332320323868323, 0.2772473834665092, 0.2338684423535113, 2.274067798290787),
    @(1.265204467861906, 0.3714242025238832, 0.03153428999343788, 0.6694046255334882, 0.002415283420845941, 0.6100265197959409, 0.2805008233742328, 0.2752309631595864, 2.246946597948408),
    @(1.464985878237485, 0.4036288395750773, 0.03496468616783943, 0.6740629386099144, 0.002411023081879477, 0.5952853035930303, 0.2836599993357112, 0.3058351721994157, 2.236034207865742),
    @(1.555753009756984, 0.4182183955511505, 0.03651768614905393, 0.6767668892621188, 0.002409178236891074, 0.5890823494178186, 0.2852641352805421, 0.3198030320214116, 2.23303850722823),
    @(1.590106281127476, 0.4237341037578233, 0.03710466097014375, 0.6778751534139218, 0.002408492974429943, 0.5868059511484738, 0.2858955802616947, 0.3250986949147716, 2.232188041632895),
    @(1.582708530860202, 0.4225466047587929, 0.03697829545125586, 0.6776327142907235, 0.002408639965612561, 0.587292985307748, 0.2857585205274091, 0.3239579013513847, 2.23235856217579),
    @(1.558579650602837, 0.4186723593010697, 0.03656599940901373, 0.6768563752521359, 0.002409121592844565, 0.5888936140376302, 0.2853156039018074, 0.3202385838407977, 2.232962843241751),
    @(1.543797581618605, 0.4162980847517872, 0.03631331030119611, 0.6763918352267524, 0.002409418338863919, 0.5898834965690654, 0.2850474286757816, 0.3179612122409097, 2.233369986623444),
    @(1.459051574445198, 0.4026741312545425, 0.03486304009493324, 0.6738980180611662, 0.002411145517221387, 0.5957008187802266, 0.2835585240486296, 0.3049232404474012, 2.236269669196929),
    @(1.407032037789236, 0.3943005393084604, 0.03397140035750823, 0.6725181063844445, 0.002412228912529685, 0.5993985113554423, 0.2826878888421334, 0.2969364281120761, 2.238553318387943),
    @(1.377101112682169, 0.3894785968053327, 0.03345784839261512, 0.6717794536598589, 0.002412860829232732, 0.6015726490642557, 0.2822028437100386, 0.2923469558539722, 2.240052035424014),
    @(1.366965239839146, 0.387845005868428, 0.03328384856772715, 0.6715388038422603, 0.002413076294739203, 0.6023168960679506, 0.2820413165782725, 0.2907937910773555, 2.240591260263244),
    @(1.412570723105375, 0.3951925121728834, 0.03406639007225465, 0.6726593026150738, 0.002412112675293694, 0.5989999864208286, 0.2827789425229241, 0.2977861915419098, 2.238291043756107),
    @(1.565667397906736, 0.4198105671508472, 0.03668713121060563, 0.6770821141142704, 0.002408979765690943, 0.5884215005815499, 0.2854450482535213, 0.3213308675831854, 2.232777638114015),
    @(1.665617593062052, 0.4358470411046937, 0.03839343681116247, 0.6804643633294916, 0.00240700995751789, 0.5819306815948764, 0.2873273445532334, 0.3367554968867381, 2.230829743679863),
    @(1.61228270987084, 0.4272930205909518, 0.03748335515483348, 0.6786141270438861, 0.002408054189323661, 0.585356195500438, 0.2863099394734974, 0.3285197967636861, 2.231717607305541),
    @(1.4100667587212, 0.3947892756775957, 0.03402344813654423, 0.6725952975017151, 0.002412165198006074, 0.599180009153276, 0.2827377288886765, 0.2974020064601532, 2.238409039372897),
    @(1.191543115048432, 0.359513082455976, 0.03026457303946728, 0.6682061165370357, 0.002416935396684027, 0.6158996401714809, 0.2794856391104332, 0.2752309631595864, 2.252704503005617)
)

$cols = @('B','C','D','F','G','I','L','M','O')

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $i + 2
    $rowVals = $data[$i]
    for ($j = 0; $j -lt $cols.Count; $j++) {
        $cellAddr = "{0}{1}" -f $cols[$j], $row
        $ws.Range($cellAddr).Value = $rowVals[$j]
    }
}
